$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.35379
$ws.Cells.Item(2, 8).Value = 1.06137
$ws.Cells.Item(2, 9).Value = 0.008539219707589354
$ws.Cells.Item(2, 10).Value = 0.008539219707589352
$ws.Cells.Item(2, 13).Value = 46.63275166666667
$ws.Cells.Item(2, 14).Value = 139.898255
$ws.Cells.Item(2, 15).Value = 0.9158911059585902
$ws.Cells.Item(2, 16).Value = 0.9158911059585902
$ws.Cells.Item(2, 17).Value = 16.49820121215
$ws.Cells.Item(2, 18).Value = 148.48381090935
$ws.Cells.Item(2, 19).Value = 0.007820995382007404
$ws.Cells.Item(2, 20).Value = 0.007820995382007402
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.35379
$ws.Cells.Item(3, 8).Value = 1.06137
$ws.Cells.Item(3, 9).Value = 0.008539219707589354
$ws.Cells.Item(3, 10).Value = 0.008539219707589352
$ws.Cells.Item(3, 15).Value = 0.05441917700612491
$ws.Cells.Item(3, 16).Value = 0.05441917700612491
$ws.Cells.Item(3, 17).Value = 0.9802677700500001
$ws.Cells.Item(3, 18).Value = 8.82240993045
$ws.Cells.Item(3, 19).Value = 0.0004646973087614952
$ws.Cells.Item(3, 20).Value = 0.0004646973087614952
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.35379
$ws.Cells.Item(4, 8).Value = 1.06137
$ws.Cells.Item(4, 9).Value = 0.008539219707589354
$ws.Cells.Item(4, 10).Value = 0.008539219707589352
$ws.Cells.Item(4, 13).Value = 0.8496050000000001
$ws.Cells.Item(4, 14).Value = 2.548815
$ws.Cells.Item(4, 15).Value = 0.01668667696558362
$ws.Cells.Item(4, 16).Value = 0.01668667696558362
$ws.Cells.Item(4, 17).Value = 0.3005817529500001
$ws.Cells.Item(4, 18).Value = 2.705235776550001
$ws.Cells.Item(4, 19).Value = 0.000142491200798689
$ws.Cells.Item(4, 20).Value = 0.000142491200798689
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.35379
$ws.Cells.Item(5, 8).Value = 1.06137
$ws.Cells.Item(5, 9).Value = 0.008539219707589354
$ws.Cells.Item(5, 10).Value = 0.008539219707589352
$ws.Cells.Item(5, 13).Value = 0.662052
$ws.Cells.Item(5, 14).Value = 1.986156
$ws.Cells.Item(5, 15).Value = 0.01300304006970129
$ws.Cells.Item(5, 16).Value = 0.0130030400697013
$ws.Cells.Item(5, 17).Value = 0.23422737708
$ws.Cells.Item(5, 18).Value = 2.10804639372
$ws.Cells.Item(5, 19).Value = 0.0001110358160217673
$ws.Cells.Item(5, 20).Value = 0.0001110358160217673
$ws.Cells.Item(6, 7).Value = 36.97491766666666
$ws.Cells.Item(6, 9).Value = 0.8924416903408624
$ws.Cells.Item(6, 10).Value = 0.8924416903408623
$ws.Cells.Item(6, 13).Value = 46.63275166666667
$ws.Cells.Item(6, 14).Value = 139.898255
$ws.Cells.Item(6, 15).Value = 0.9158911059585902
$ws.Cells.Item(6, 16).Value = 0.9158911059585902
$ws.Cells.Item(6, 17).Value = 1724.242153445113
$ws.Cells.Item(6, 18).Value = 15518.17938100602
$ws.Cells.Item(6, 19).Value = 0.8173794067698462
$ws.Cells.Item(6, 20).Value = 0.8173794067698461
$ws.Cells.Item(7, 7).Value = 36.97491766666666
$ws.Cells.Item(7, 9).Value = 0.8924416903408624
$ws.Cells.Item(7, 10).Value = 0.8924416903408623
$ws.Cells.Item(7, 15).Value = 0.05441917700612491
$ws.Cells.Item(7, 16).Value = 0.05441917700612491
$ws.Cells.Item(7, 18).Value = 922.0381604906049
$ws.Cells.Item(7, 19).Value = 0.0485659423143047
$ws.Cells.Item(7, 20).Value = 0.0485659423143047
$ws.Cells.Item(8, 7).Value = 36.97491766666666
$ws.Cells.Item(8, 9).Value = 0.8924416903408624
$ws.Cells.Item(8, 10).Value = 0.8924416903408623
$ws.Cells.Item(8, 13).Value = 0.8496050000000001
$ws.Cells.Item(8, 14).Value = 2.548815
$ws.Cells.Item(8, 15).Value = 0.01668667696558362
$ws.Cells.Item(8, 16).Value = 0.01668667696558362
$ws.Cells.Item(8, 17).Value = 31.41407492418833
$ws.Cells.Item(8, 18).Value = 282.726674317695
$ws.Cells.Item(8, 19).Value = 0.01489188619733738
$ws.Cells.Item(8, 20).Value = 0.01489188619733738
$ws.Cells.Item(9, 7).Value = 36.97491766666666
$ws.Cells.Item(9, 9).Value = 0.8924416903408624
$ws.Cells.Item(9, 10).Value = 0.8924416903408623
$ws.Cells.Item(9, 13).Value = 0.662052
$ws.Cells.Item(9, 14).Value = 1.986156
$ws.Cells.Item(9, 15).Value = 0.01300304006970129
$ws.Cells.Item(9, 16).Value = 0.0130030400697013
$ws.Cells.Item(9, 17).Value = 24.47931819105199
$ws.Cells.Item(9, 18).Value = 220.313863719468
$ws.Cells.Item(9, 19).Value = 0.01160445505937419
$ws.Cells.Item(9, 20).Value = 0.01160445505937419
$ws.Cells.Item(10, 7).Value = 4.102478333333333
$ws.Cells.Item(10, 8).Value = 12.307435
$ws.Cells.Item(10, 9).Value = 0.09901908995154843
$ws.Cells.Item(10, 10).Value = 0.0990190899515484
$ws.Cells.Item(10, 13).Value = 46.63275166666667
$ws.Cells.Item(10, 14).Value = 139.898255
$ws.Cells.Item(10, 15).Value = 0.9158911059585902
$ws.Cells.Item(10, 16).Value = 0.9158911059585902
$ws.Cells.Item(10, 17).Value = 191.3098533362139
$ws.Cells.Item(10, 18).Value = 1721.788680025925
$ws.Cells.Item(10, 19).Value = 0.09069070380673681
$ws.Cells.Item(10, 20).Value = 0.0906907038067368
$ws.Cells.Item(11, 7).Value = 4.102478333333333
$ws.Cells.Item(11, 8).Value = 12.307435
$ws.Cells.Item(11, 9).Value = 0.09901908995154843
$ws.Cells.Item(11, 10).Value = 0.0990190899515484
$ws.Cells.Item(11, 15).Value = 0.05441917700612491
$ws.Cells.Item(11, 16).Value = 0.05441917700612491
$ws.Cells.Item(11, 17).Value = 11.36698970433055
$ws.Cells.Item(11, 18).Value = 102.302907338975
$ws.Cells.Item(11, 19).Value = 0.005388537383058718
$ws.Cells.Item(11, 20).Value = 0.005388537383058717
$ws.Cells.Item(12, 7).Value = 4.102478333333333
$ws.Cells.Item(12, 8).Value = 12.307435
$ws.Cells.Item(12, 9).Value = 0.09901908995154843
$ws.Cells.Item(12, 10).Value = 0.0990190899515484
$ws.Cells.Item(12, 13).Value = 0.8496050000000001
$ws.Cells.Item(12, 14).Value = 2.548815
$ws.Cells.Item(12, 15).Value = 0.01668667696558362
$ws.Cells.Item(12, 16).Value = 0.01668667696558362
$ws.Cells.Item(12, 17).Value = 3.485486104391667
$ws.Cells.Item(12, 18).Value = 31.369374939525
$ws.Cells.Item(12, 19).Value = 0.001652299567447556
$ws.Cells.Item(12, 20).Value = 0.001652299567447556
$ws.Cells.Item(13, 7).Value = 4.102478333333333
$ws.Cells.Item(13, 8).Value = 12.307435
$ws.Cells.Item(13, 9).Value = 0.09901908995154843
$ws.Cells.Item(13, 10).Value = 0.0990190899515484
$ws.Cells.Item(13, 13).Value = 0.662052
$ws.Cells.Item(13, 14).Value = 1.986156
$ws.Cells.Item(13, 15).Value = 0.01300304006970129
$ws.Cells.Item(13, 16).Value = 0.0130030400697013
$ws.Cells.Item(13, 17).Value = 2.71605398554
$ws.Cells.Item(13, 18).Value = 24.44448586986
$ws.Cells.Item(13, 19).Value = 0.001287549194305341
$ws.Cells.Item(13, 20).Value = 0.001287549194305341

Write-Output "Updated cells: 141"